$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column B values for rows 2 through 73 (naive forecaster bugfix)
$ws.Cells.Item(2, 2).Value = 0.25
$ws.Cells.Item(3, 2).Value = 0.25
$ws.Cells.Item(4, 2).Value = 0.1500000000000057
$ws.Cells.Item(5, 2).Value = 0.2000000000000028
$ws.Cells.Item(6, 2).Value = 0.4000000000000057
$ws.Cells.Item(7, 2).Value = 0.4999999999999858
$ws.Cells.Item(8, 2).Value = 0.09999999999999432
$ws.Cells.Item(9, 2).Value = 0.4000000000000057
$ws.Cells.Item(10, 2).Value = -0.5
$ws.Cells.Item(11, 2).Value = 0.4999999999999858
$ws.Cells.Item(12, 2).Value = 0.2999999999999829
$ws.Cells.Item(13, 2).Value = 0.2250000000000085
$ws.Cells.Item(14, 2).Value = 0.4
$ws.Cells.Item(15, 2).Value = 0.2
$ws.Cells.Item(16, 2).Value = 0.2
$ws.Cells.Item(17, 2).Value = 0.3
$ws.Cells.Item(18, 2).Value = 0.3
$ws.Cells.Item(19, 2).Value = 0.3
$ws.Cells.Item(20, 2).Value = 0.4
$ws.Cells.Item(21, 2).Value = 0.3
$ws.Cells.Item(22, 2).Value = 0.1
$ws.Cells.Item(23, 2).Value = 0.2
$ws.Cells.Item(24, 2).Value = 0.2
$ws.Cells.Item(25, 2).Value = 0.3
$ws.Cells.Item(26, 2).Value = 0.3
$ws.Cells.Item(27, 2).Value = 0.4
$ws.Cells.Item(28, 2).Value = 0.3490000000000038
$ws.Cells.Item(29, 2).Value = 0.4399999999999977
$ws.Cells.Item(30, 2).Value = 0.4399999999999977
$ws.Cells.Item(31, 2).Value = 0.3500000000000085
$ws.Cells.Item(32, 2).Value = 0.4499999999999886
$ws.Cells.Item(33, 2).Value = 0.4000000000000057
$ws.Cells.Item(34, 2).Value = 0.4202440737484352
$ws.Cells.Item(35, 2).Value = 0.4753723183093825
$ws.Cells.Item(36, 2).Value = 0.2999999999999687
$ws.Cells.Item(37, 2).Value = 0.4499999999999886
$ws.Cells.Item(38, 2).Value = 0.4000000000000057
$ws.Cells.Item(39, 2).Value = 0.3500000000000085
$ws.Cells.Item(40, 2).Value = 0.4499999999999744
$ws.Cells.Item(41, 2).Value = 0.4250000000000114
$ws.Cells.Item(42, 2).Value = 0.4000000000000057
$ws.Cells.Item(43, 2).Value = 0.4000000000000057
$ws.Cells.Item(44, 2).Value = 0.4000000000000057
$ws.Cells.Item(45, 2).Value = 0.4000000000000057
$ws.Cells.Item(46, 2).Value = 0.3
$ws.Cells.Item(47, 2).Value = 0.539999999999992
$ws.Cells.Item(48, 2).Value = -0.3
$ws.Cells.Item(49, 2).Value = 0.0999999999999659
$ws.Cells.Item(50, 2).Value = 0.1399999999999864
$ws.Cells.Item(51, 2).Value = -0.0999999999999659
$ws.Cells.Item(52, 2).Value = -10.8
$ws.Cells.Item(53, 2).Value = 7.25
$ws.Cells.Item(54, 2).Value = -0.7999999999999972
$ws.Cells.Item(55, 2).Value = -2.400000000000006
$ws.Cells.Item(56, 2).Value = 3.480874220397794
$ws.Cells.Item(57, 2).Value = 6.400000000000006
$ws.Cells.Item(58, 2).Value = -1.099999999999994
$ws.Cells.Item(59, 2).Value = 0.4209467346675666
$ws.Cells.Item(60, 2).Value = 1.540000000000006
$ws.Cells.Item(61, 2).Value = -0.1490000000000009
$ws.Cells.Item(62, 2).Value = -0.4999966213670604
$ws.Cells.Item(63, 2).Value = -0.7399999999999949
$ws.Cells.Item(64, 2).Value = -0.2510000000000048
$ws.Cells.Item(65, 2).Value = 0.2000000000000028
$ws.Cells.Item(66, 2).Value = 0.1200000000000045
$ws.Cells.Item(67, 2).Value = 0.09999999999999432
$ws.Cells.Item(68, 2).Value = 0
$ws.Cells.Item(69, 2).Value = 0.09999999999999432
$ws.Cells.Item(70, 2).Value = 0.3200015876295765
$ws.Cells.Item(71, 2).Value = 0.09999771493470178
$ws.Cells.Item(72, 2).Value = 0.2000034419242951
$ws.Cells.Item(73, 2).Value = 0.03999999999999204

# Remove now-obsolete trailing rows 74 to 82
$ws.Range("A74:B82").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp) | Out-Null

$wb.Save()
